$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second part of project -> uncertainty analysis
$ws.Range("F5").Value = "no"
$ws.Range("E11").Value = 12.2533526453471
$ws.Range("F11").Value = "no"

# Update active selection
$ws.Range("J9").Select()
